$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 68

# Column A stores dates as literal text (e.g. "09/03/2025"), not real Excel
# dates. Typing a date-like string directly into a General-formatted cell
# would make Excel auto-convert it into a date serial number, so instead
# we write it as a formula that evaluates to the text string, then paste
# the result back as a plain value. This keeps the cell a literal string
# (matching every other row in the column) without leaving behind any new
# cell style.
$ws.Cells.Item($row, 1).Formula = "=""11/08/2025"""
$ws.Cells.Item($row, 1).Copy()
$ws.Cells.Item($row, 1).PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = 0

$ws.Cells.Item($row, 2).Value = 0.185766309476215
$ws.Cells.Item($row, 3).Value = 0.814233690523785
